$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17 (hunk 0)
$ws.Range("H17").Value = 960
$ws.Range("J17").Value = 1011.1111
$ws.Range("L17").Value = 3033.3333
$ws.Range("N17").Value = -3369.3333

# Row 28 (hunk 1)
$ws.Range("H28").Value = 7774.375
$ws.Range("I28").Value = 7580.2
$ws.Range("K28").Value = 7580.2
$ws.Range("M28").Value = -7095.2

# Row 32 (hunk 2)
$ws.Range("H32").Value = 7034.4443
$ws.Range("I32").Value = 466.66666
$ws.Range("J32").Value = 10318.333
$ws.Range("K32").Value = 466.66666
$ws.Range("L32").Value = 10318.333
$ws.Range("M32").Value = -140.66666
$ws.Range("N32").Value = -10970.333

# Row 38 (hunk 3)
$ws.Range("H38").Value = 278.41666
$ws.Range("I38").Value = 122
$ws.Range("K38").Value = 366
$ws.Range("M38").Value = 6

# Row 70 (hunk 4)
$ws.Range("H70").Value = 1037.375
$ws.Range("I70").Value = 800.5
$ws.Range("J70").Value = 1071.2142
$ws.Range("K70").Value = 2401.5
$ws.Range("L70").Value = 3213.6426
$ws.Range("M70").Value = -2131.5
$ws.Range("N70").Value = -3753.6426

# Row 73 (hunk 5)
$ws.Range("H73").Value = 1037.375
$ws.Range("I73").Value = 800.5
$ws.Range("J73").Value = 1071.2142
$ws.Range("K73").Value = 2401.5
$ws.Range("L73").Value = 3213.6426
$ws.Range("M73").Value = -1465.5
$ws.Range("N73").Value = -5085.642599999999

# Row 113 (hunk 6)
$ws.Range("H113").Value = 2483.1875
$ws.Range("I113").Value = 2447.3076
$ws.Range("J113").Value = 2638.6667
$ws.Range("K113").Value = 2447.3076
$ws.Range("L113").Value = 2638.6667
$ws.Range("M113").Value = 806.6923999999999
$ws.Range("N113").Value = -9146.6667

# Row 132 (hunk 7)
$ws.Range("H132").Value = 9531333
$ws.Range("I132").Value = 10758701
$ws.Range("J132").Value = 19226.5
$ws.Range("K132").Value = 32276103
$ws.Range("L132").Value = 57679.5
$ws.Range("M132").Value = -32273573
$ws.Range("N132").Value = -62739.5

# Row 137 (hunk 8)
$ws.Range("H137").Value = 1962.6562
$ws.Range("J137").Value = 3157.8333
$ws.Range("L137").Value = 9473.499899999999
$ws.Range("N137").Value = -14573.4999

# Row 138 (hunk 9)
$ws.Range("H138").Value = 2720.8706
$ws.Range("I138").Value = 2055.5
$ws.Range("J138").Value = 2852.0703
$ws.Range("K138").Value = 6166.5
$ws.Range("L138").Value = 8556.2109
$ws.Range("M138").Value = -1026.5
$ws.Range("N138").Value = -18836.2109

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 10)
$ws.Range("H32").Value = 13167.036
$ws.Range("I32").Value = 10108.9795
$ws.Range("J32").Value = 34573.43
$ws.Range("K32").Value = 10108.9795
$ws.Range("L32").Value = 34573.43
$ws.Range("M32").Value = -9821.979499999999
$ws.Range("N32").Value = -35147.43

# Row 45 (hunk 11)
$ws.Range("H45").Value = 1224.3125
$ws.Range("I45").Value = 989.9091
$ws.Range("J45").Value = 1740
$ws.Range("K45").Value = 989.9091
$ws.Range("L45").Value = 1740
$ws.Range("M45").Value = -612.9091
$ws.Range("N45").Value = -2494

# Row 61 (hunk 12)
$ws.Range("H61").Value = 45455756
$ws.Range("I61").Value = 58824388
$ws.Range("J61").Value = 2400
$ws.Range("K61").Value = 58824388
$ws.Range("L61").Value = 2400
$ws.Range("M61").Value = -58824176
$ws.Range("N61").Value = -2824

# Row 102 (hunk 13)
$ws.Range("H102").Value = 23811840
$ws.Range("I102").Value = 33334976
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 33334976
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -33333354
$ws.Range("N102").Value = -7244

# Row 110 (hunk 14)
$ws.Range("H110").Value = 1056
$ws.Range("I110").Value = 1100
$ws.Range("J110").Value = 990
$ws.Range("K110").Value = 1100
$ws.Range("L110").Value = 990
$ws.Range("M110").Value = 945
$ws.Range("N110").Value = -5080

# Row 136 (hunk 15)
$ws.Range("H136").Value = 45455756
$ws.Range("I136").Value = 58824388
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 176473164
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -176470614
$ws.Range("N136").Value = -12300

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86 (hunk 16)
$ws.Range("H86").Value = 3780.75
$ws.Range("I86").Value = 3680.5334
$ws.Range("J86").Value = 4081.4
$ws.Range("K86").Value = 3680.5334
$ws.Range("L86").Value = 4081.4
$ws.Range("M86").Value = -2557.5334
$ws.Range("N86").Value = -6327.4

# Row 89 (hunk 17)
$ws.Range("H89").Value = 3780.75
$ws.Range("I89").Value = 3680.5334
$ws.Range("J89").Value = 4081.4
$ws.Range("K89").Value = 18402.667
$ws.Range("L89").Value = 20407
$ws.Range("M89").Value = -12786.667
$ws.Range("N89").Value = -31639

# Row 105 (hunk 18)
$ws.Range("H105").Value = 144270700
$ws.Range("I105").Value = 168315490
$ws.Range("K105").Value = 168315490
$ws.Range("M105").Value = -168313743

# Row 107 (hunk 19)
$ws.Range("H107").Value = 825.8095
$ws.Range("I107").Value = 836.1875
$ws.Range("K107").Value = 836.1875
$ws.Range("M107").Value = 1083.8125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (hunk 20)
$ws.Range("H31").Value = 1498.7693
$ws.Range("I31").Value = 1498.7693
$ws.Range("K31").Value = 1498.7693
$ws.Range("M31").Value = -1203.7693

# Row 34 (hunk 21)
$ws.Range("H34").Value = 1498.7693
$ws.Range("I34").Value = 1498.7693
$ws.Range("K34").Value = 1498.7693
$ws.Range("M34").Value = -1296.7693

# Row 68 (hunk 22)
$ws.Range("H68").Value = 11750
$ws.Range("J68").Value = 11750
$ws.Range("L68").Value = 11750
$ws.Range("N68").Value = -13248

# Row 71 (hunk 23)
$ws.Range("H71").Value = 11750
$ws.Range("J71").Value = 11750
$ws.Range("L71").Value = 35250
$ws.Range("N71").Value = -42738

# Row 74 (hunk 24)
$ws.Range("H74").Value = 26666.666
$ws.Range("I74").Value = 23500
$ws.Range("K74").Value = 23500
$ws.Range("M74").Value = -22626

# Row 77 (hunk 25)
$ws.Range("H77").Value = 26666.666
$ws.Range("I77").Value = 23500
$ws.Range("K77").Value = 70500
$ws.Range("M77").Value = -66132

# Row 99 (hunk 26)
$ws.Range("H99").Value = 1699.2
$ws.Range("I99").Value = 1686.125
$ws.Range("J99").Value = 1714.1428
$ws.Range("K99").Value = 1686.125
$ws.Range("L99").Value = 1714.1428
$ws.Range("M99").Value = -188.125
$ws.Range("N99").Value = -4710.1428

# Row 107 (hunk 27)
$ws.Range("H107").Value = 1491.8
$ws.Range("I107").Value = 815.6
$ws.Range("J107").Value = 2168
$ws.Range("K107").Value = 815.6
$ws.Range("L107").Value = 2168
$ws.Range("M107").Value = 1104.4
$ws.Range("N107").Value = -6008

# Row 109 (hunk 28)
$ws.Range("H109").Value = 14667
$ws.Range("J109").Value = 14667
$ws.Range("L109").Value = 14667
$ws.Range("N109").Value = -16747

# Row 126 (hunk 29)
$ws.Range("H126").Value = 1699.2
$ws.Range("I126").Value = 1686.125
$ws.Range("J126").Value = 1714.1428
$ws.Range("K126").Value = 5058.375
$ws.Range("L126").Value = 5142.428400000001
$ws.Range("M126").Value = -2588.375
$ws.Range("N126").Value = -10082.4284

# Row 132 (hunk 30)
$ws.Range("H132").Value = 1715.3948
$ws.Range("I132").Value = 1357.375
$ws.Range("K132").Value = 4072.125
$ws.Range("M132").Value = -1542.125

# Row 134 (hunk 31)
$ws.Range("H134").Value = 14707564
$ws.Range("I134").Value = 1753.1904
$ws.Range("J134").Value = 38463104
$ws.Range("K134").Value = 5259.5712
$ws.Range("L134").Value = 115389312
$ws.Range("M134").Value = -2724.5712
$ws.Range("N134").Value = -115394382

# Row 135 (hunk 32)
$ws.Range("H135").Value = 37706.25
$ws.Range("J135").Value = 37706.25
$ws.Range("L135").Value = 37706.25
$ws.Range("N135").Value = -47846.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 107 (hunk 33)
$ws.Range("H107").Value = 7711.9287
$ws.Range("I107").Value = 556.6
$ws.Range("J107").Value = 11687.111
$ws.Range("K107").Value = 1669.8
$ws.Range("L107").Value = 35061.333
$ws.Range("M107").Value = 250.1999999999998
$ws.Range("N107").Value = -38901.333

# Row 109 (hunk 34)
$ws.Range("H109").Value = 169337.83
$ws.Range("I109").Value = 501013.5
$ws.Range("K109").Value = 1503040.5
$ws.Range("M109").Value = -1502000.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113 (hunk 35)
$ws.Range("H113").Value = 1489
$ws.Range("I113").Value = 1298
$ws.Range("J113").Value = 1680
$ws.Range("K113").Value = 1298
$ws.Range("L113").Value = 1680
$ws.Range("M113").Value = 872
$ws.Range("N113").Value = -6020

# Row 126 (hunk 36)
$ws.Range("H126").Value = 2118.3333
$ws.Range("I126").Value = 1952.5
$ws.Range("J126").Value = 2450
$ws.Range("K126").Value = 5857.5
$ws.Range("L126").Value = 7350
$ws.Range("M126").Value = -3387.5
$ws.Range("N126").Value = -12290

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (hunk 37)
$ws.Range("H7").Value = 2579.8333
$ws.Range("I7").Value = 1993.3334
$ws.Range("J7").Value = 3166.3333
$ws.Range("K7").Value = 1993.3334
$ws.Range("L7").Value = 3166.3333
$ws.Range("M7").Value = -1881.3334
$ws.Range("N7").Value = -3390.3333

# Row 10 (hunk 38)
$ws.Range("H10").Value = 100
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# Row 22 (hunk 39)
$ws.Range("H22").Value = 999.6667
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1249.5
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 1249.5
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -1839.5

# Row 27 (hunk 40)
$ws.Range("H27").Value = 999.6667
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 1249.5
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 1249.5
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -1463.5

# Row 80 (hunk 41)
$ws.Range("H80").Value = 20128
$ws.Range("J80").Value = 20128
$ws.Range("L80").Value = 20128
$ws.Range("N80").Value = -22374

# Row 83 (hunk 42)
$ws.Range("H83").Value = 20128
$ws.Range("J83").Value = 20128
$ws.Range("L83").Value = 60384
$ws.Range("N83").Value = -71616

# Row 126 (hunk 43)
$ws.Range("H126").Value = 2579.8333
$ws.Range("I126").Value = 1993.3334
$ws.Range("J126").Value = 3166.3333
$ws.Range("K126").Value = 5980.0002
$ws.Range("L126").Value = 9498.999899999999
$ws.Range("M126").Value = -3510.0002
$ws.Range("N126").Value = -14438.9999

# Row 132 (hunk 44)
$ws.Range("H132").Value = 2240.775
$ws.Range("I132").Value = 1820.5
$ws.Range("J132").Value = 3021.2856
$ws.Range("K132").Value = 5461.5
$ws.Range("L132").Value = 9063.856800000001
$ws.Range("M132").Value = -2931.5
$ws.Range("N132").Value = -14123.8568

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 11 (hunk 45)
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# Row 41 (hunk 46)
$ws.Range("H41").Value = 10399.5
$ws.Range("J41").Value = 10399.5
$ws.Range("L41").Value = 10399.5
$ws.Range("N41").Value = -11179.5

# Row 107 (hunk 47)
$ws.Range("H107").Value = 500350
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 1000000
$ws.Range("K107").Value = 2100
$ws.Range("L107").Value = 3000000
$ws.Range("M107").Value = -180
$ws.Range("N107").Value = -3003840

# Row 122 (hunk 48)
$ws.Range("H122").Value = 13891597
$ws.Range("I122").Value = 14708574
$ws.Range("K122").Value = 44125722
$ws.Range("M122").Value = -44123272
